# data_import now skips procedureItems whose text is only whitespace.
# For the "Analysephase" method row (German sheet) that means the two
# previously-skipped procedure items - "teilnehmende Beobachtungen" and
# "Personas" - are now imported as two additional columns, and every
# existing Analysephase row gets a value recorded for them.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("German")

# --- new header columns (row 1) ------------------------------------------
$ws.Range("R1").Value2 = "Beobachtungen der Teilnehmenden"
$ws.Range("S1").Value2 = "Personas"

# --- row 2 (Einzelinterviews / id 0) --------------------------------------
$ws.Range("R2").Value2 = "Teilnehmende Beobachtungen"
$ws.Range("S2").Value2 = "Personas"

# --- copy the "Entscheidender Erfolgsfaktor..." description (with its
#     wrap-text formatting) from C2 down into C3 and C4 -------------------
$ws.Range("C2").AutoFill($ws.Range("C2:C4"), 0)

# --- row 3 (Methoden fuer die Konzeptionsphase / id 1) --------------------
$ws.Range("R3").Value2 = "Teilnehmende Beobachtungen"
$ws.Range("S3").Value2 = "Personas"

# --- row 4 (Umsetzungs- und Evaluationsphase / id 2) -----------------------
$ws.Range("R4").Value2 = "Teilnehmende Beobachtungen"
$ws.Range("S4").Value2 = "Personas"

# --- move the view so the newly added columns are visible ------------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("S2").Select()

Write-Host "German sheet updated: columns R:S added through row 4"
